# updated CB_API and Dash
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2 updates
$ws.Range("S2").Value = 0.01015625
$ws.Range("T2").Value = 0.1911458333333333
$ws.Range("U2").Value = 1.3353125
$ws.Range("V2").Value = 6.1259375

# Row 14 updates
$ws.Range("J14").Value = 7.605744791666667
$ws.Range("K14").Value = 11.4377890625
$ws.Range("L14").Value = 24.758015625
$ws.Range("M14").Value = 18.6651171875
$ws.Range("N14").Value = 5.311622395833334
$ws.Range("O14").Value = 6.1869375
$ws.Range("P14").Value = 10.64615885416667
$ws.Range("Q14").Value = 10.20610416666667
$ws.Range("R14").Value = 17.3760234375
$ws.Range("S14").Value = 22.0525
$ws.Range("T14").Value = 38.00716145833334
$ws.Range("U14").Value = 17.11038020833334
$ws.Range("V14").Value = 9.284989583333333
$ws.Range("W14").Value = 1.228276041666667
